$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1636
$ws.Range("C3").Value = 1411
$ws.Range("C4").Value = 689
$ws.Range("C5").Value = 1471
$ws.Range("C6").Value = 598
